# redo FR pop dens calcs using overlays
$wb = $excel.ActiveWorkbook

# --- Sheet "area_mixre": update descriptive stats for area column ---
$ws1 = $wb.Worksheets.Item("area_mixre")
$ws1.Range("B2").Value = 180
$ws1.Range("B3").Value = 3.05134641844613
$ws1.Range("B4").Value = 2.853763647512013
$ws1.Range("B5").Value = 0.1686295793966096
$ws1.Range("B6").Value = 0.9002816795470007
$ws1.Range("B7").Value = 1.955709873699787
$ws1.Range("B8").Value = 4.723393760657881

# --- Sheet "area_hires": update descriptive stats for area column ---
$ws2 = $wb.Worksheets.Item("area_hires")
$ws2.Range("B3").Value = 0.8930839327557512
$ws2.Range("B4").Value = 1.687649362145632
$ws2.Range("B5").Value = 0.00007841371226228017

# --- Sheet "area_pop_sum": update population & density totals ---
$ws4 = $wb.Worksheets.Item("area_pop_sum")
$ws4.Range("B3").Value = 1357466
$ws4.Range("B4").Value = 2471.526187461885
